$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 275.72726
$ws.Range("I4").Value = 209.4
$ws.Range("K4").Value = 209.4
$ws.Range("M4").Value = -95.40000000000001

$ws.Range("H9").Value = 8291.076999999999
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

$ws.Range("H12").Value = 16914.5
$ws.Range("I12").Value = 20198.4
$ws.Range("K12").Value = 20198.4
$ws.Range("M12").Value = -20028.4

$ws.Range("H33").Value = 285.33334
$ws.Range("I33").Value = 303.9
$ws.Range("K33").Value = 303.9
$ws.Range("M33").Value = -74.89999999999998

$ws.Range("H40").Value = 12349234
$ws.Range("I40").Value = 1999.6
$ws.Range("J40").Value = 27783278
$ws.Range("K40").Value = 1999.6
$ws.Range("L40").Value = 27783278
$ws.Range("M40").Value = -1824.6
$ws.Range("N40").Value = -27783628

$ws.Range("H92").Value = 866.5833
$ws.Range("I92").Value = 854.4545000000001
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 854.4545000000001
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 393.5454999999999
$ws.Range("N92").Value = -3496

$ws.Range("H99").Value = 1497.5
$ws.Range("I99").Value = 460.75
$ws.Range("J99").Value = 5644.5
$ws.Range("K99").Value = 1382.25
$ws.Range("L99").Value = 16933.5
$ws.Range("M99").Value = 115.75
$ws.Range("N99").Value = -19929.5

$ws.Range("H103").Value = 1437.6
$ws.Range("I103").Value = 1268.4286
$ws.Range("J103").Value = 1832.3334
$ws.Range("K103").Value = 3805.2858
$ws.Range("L103").Value = 5497.0002
$ws.Range("M103").Value = -3219.2858
$ws.Range("N103").Value = -6669.0002

$ws.Range("H112").Value = 74856.71000000001
$ws.Range("I112").Value = 2333.3333
$ws.Range("J112").Value = 94635.82000000001
$ws.Range("K112").Value = 6999.999899999999
$ws.Range("L112").Value = 283907.46
$ws.Range("M112").Value = -5891.999899999999
$ws.Range("N112").Value = -286123.46

$ws.Range("H116").Value = 10013.75
$ws.Range("I116").Value = 10414.267
$ws.Range("J116").Value = 4006
$ws.Range("K116").Value = 10414.267
$ws.Range("L116").Value = 4006
$ws.Range("M116").Value = -6972.267
$ws.Range("N116").Value = -10890

$ws.Range("H137").Value = 2025.9048
$ws.Range("I137").Value = 1408.5294
$ws.Range("J137").Value = 4649.75
$ws.Range("K137").Value = 4225.5882
$ws.Range("L137").Value = 13949.25
$ws.Range("M137").Value = -1675.5882
$ws.Range("N137").Value = -19049.25

$ws.Range("H138").Value = 2502.2273
$ws.Range("I138").Value = 1870.3334
$ws.Range("J138").Value = 3079.1738
$ws.Range("K138").Value = 5611.0002
$ws.Range("L138").Value = 9237.5214
$ws.Range("M138").Value = -471.0002000000004
$ws.Range("N138").Value = -19517.5214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 12681.75
$ws.Range("I4").Value = 135.83333
$ws.Range("K4").Value = 135.83333
$ws.Range("M4").Value = -19.83332999999999

$ws.Range("H6").Value = 45334
$ws.Range("I6").Value = 18001
$ws.Range("K6").Value = 18001
$ws.Range("M6").Value = -17828

$ws.Range("H45").Value = 1495
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -623

$ws.Range("H132").Value = 3132404.2
$ws.Range("I132").Value = 5270539.5
$ws.Range("K132").Value = 15811618.5
$ws.Range("M132").Value = -15809088.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6804392.5
$ws.Range("I22").Value = 299.75
$ws.Range("K22").Value = 299.75
$ws.Range("M22").Value = -126.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 86507.21000000001
$ws.Range("I25").Value = 1100
$ws.Range("J25").Value = 93077
$ws.Range("K25").Value = 1100
$ws.Range("L25").Value = 93077
$ws.Range("M25").Value = -926
$ws.Range("N25").Value = -93425

$ws.Range("H58").Value = 45467870
$ws.Range("I58").Value = 62513700
$ws.Range("K58").Value = 62513700
$ws.Range("M58").Value = -62513497

$ws.Range("H116").Value = 75000
$ws.Range("J116").Value = 75000
$ws.Range("L116").Value = 75000
$ws.Range("N116").Value = -84178

$ws.Range("H132").Value = 111112290
$ws.Range("I132").Value = 142858530
$ws.Range("J132").Value = 450
$ws.Range("K132").Value = 428575590
$ws.Range("L132").Value = 1350
$ws.Range("M132").Value = -428573060
$ws.Range("N132").Value = -6410

$ws.Range("H136").Value = 45467870
$ws.Range("I136").Value = 62513700
$ws.Range("K136").Value = 187541100
$ws.Range("M136").Value = -187538550

$ws.Range("H141").Value = 253772.73
$ws.Range("J141").Value = 381642.84
$ws.Range("L141").Value = 381642.84
$ws.Range("N141").Value = -392002.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1889128.5
$ws.Range("I4").Value = 875263.9
$ws.Range("J4").Value = 10000046
$ws.Range("K4").Value = 2625791.7
$ws.Range("L4").Value = 30000138
$ws.Range("M4").Value = -2625679.7
$ws.Range("N4").Value = -30000362

$ws.Range("H32").Value = 100001800
$ws.Range("I32").Value = 100001800
$ws.Range("K32").Value = 300005400
$ws.Range("M32").Value = -300005117

$ws.Range("H34").Value = 2097
$ws.Range("J34").Value = 4000.8572
$ws.Range("L34").Value = 12002.5716
$ws.Range("N34").Value = -12170.5716

$ws.Range("H39").Value = 4319.4
$ws.Range("I39").Value = 600
$ws.Range("J39").Value = 6799
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 20397
$ws.Range("M39").Value = -1506
$ws.Range("N39").Value = -20985

$ws.Range("H49").Value = 19000
$ws.Range("J49").Value = 19000
$ws.Range("L49").Value = 57000
$ws.Range("N49").Value = -57312

$ws.Range("H55").Value = 239.8
$ws.Range("I55").Value = 239.8
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 719.4000000000001
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -542.4000000000001
$ws.Range("N55").Value = ""

$ws.Range("H131").Value = 1822.5454
$ws.Range("J131").Value = 2354.4443
$ws.Range("L131").Value = 7063.3329
$ws.Range("N131").Value = -17143.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 16248.25
$ws.Range("I18").Value = 19999.334
$ws.Range("J18").Value = 4995
$ws.Range("K18").Value = 19999.334
$ws.Range("L18").Value = 4995
$ws.Range("M18").Value = -19706.334
$ws.Range("N18").Value = -5581

$ws.Range("H102").Value = 1845.3125
$ws.Range("J102").Value = 2649.6667
$ws.Range("L102").Value = 2649.6667
$ws.Range("N102").Value = -5893.6667

$ws.Range("H107").Value = 1223.25
$ws.Range("I107").Value = 863.8889
$ws.Range("J107").Value = 2301.3333
$ws.Range("K107").Value = 863.8889
$ws.Range("L107").Value = 2301.3333
$ws.Range("M107").Value = 1056.1111
$ws.Range("N107").Value = -6141.3333

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8600
$ws.Range("J2").Value = 8600
$ws.Range("L2").Value = 8600
$ws.Range("N2").Value = -8824

$ws.Range("H20").Value = 16699.5
$ws.Range("J20").Value = 16799.857
$ws.Range("L20").Value = 16799.857
$ws.Range("N20").Value = -17251.857

$ws.Range("H22").Value = 1898.8334
$ws.Range("I22").Value = 2256.7144
$ws.Range("J22").Value = 1397.8
$ws.Range("K22").Value = 2256.7144
$ws.Range("L22").Value = 1397.8
$ws.Range("M22").Value = -1961.7144
$ws.Range("N22").Value = -1987.8

$ws.Range("H27").Value = 1898.8334
$ws.Range("I27").Value = 2256.7144
$ws.Range("J27").Value = 1397.8
$ws.Range("K27").Value = 2256.7144
$ws.Range("L27").Value = 1397.8
$ws.Range("M27").Value = -2149.7144
$ws.Range("N27").Value = -1611.8

$ws.Range("H40").Value = 11737.9
$ws.Range("I40").Value = 7922.375
$ws.Range("K40").Value = 7922.375
$ws.Range("M40").Value = -7786.375

$ws.Range("H42").Value = 39000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 39000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 39000
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -40126

$ws.Range("H46").Value = 1291.1666
$ws.Range("I46").Value = 949.25
$ws.Range("J46").Value = 1975
$ws.Range("K46").Value = 949.25
$ws.Range("L46").Value = 1975
$ws.Range("M46").Value = -761.25
$ws.Range("N46").Value = -2351

$ws.Range("H49").Value = 39000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 39000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 39000
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -39294

$ws.Range("H55").Value = 669.4
$ws.Range("I55").Value = 257.83334
$ws.Range("J55").Value = 845.7857
$ws.Range("K55").Value = 257.83334
$ws.Range("L55").Value = 845.7857
$ws.Range("M55").Value = -84.83334000000002
$ws.Range("N55").Value = -1191.7857

$ws.Range("H76").Value = 9249.75
$ws.Range("J76").Value = 9249.75
$ws.Range("L76").Value = 9249.75
$ws.Range("N76").Value = -9925.75

$ws.Range("H79").Value = 9249.75
$ws.Range("J79").Value = 9249.75
$ws.Range("L79").Value = 9249.75
$ws.Range("N79").Value = -11589.75

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 39930
$ws.Range("J3").Value = 39930
$ws.Range("L3").Value = 39930
$ws.Range("N3").Value = -40158

$ws.Range("H39").Value = 65024.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 65024.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 65024.5
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -65850.5

$ws.Range("H122").Value = 4222.885
$ws.Range("I122").Value = 3937.0908
$ws.Range("J122").Value = 5794.75
$ws.Range("K122").Value = 11811.2724
$ws.Range("L122").Value = 17384.25
$ws.Range("M122").Value = -9361.2724
$ws.Range("N122").Value = -22284.25

$ws.Range("H132").Value = 13894399
$ws.Range("I132").Value = 16672790
$ws.Range("K132").Value = 50018370
$ws.Range("M132").Value = -50015840
